$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the new "darker gray trays" sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "frame3"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws2.Name = "frame14(dark)"

# --- Remove the leftover chart-tracking defined names (_xlchart.v1.*) ---
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
  $wb.Names.Item($i).Delete()
}

# --- Header row for the new sheet (reuses existing shared strings) ---
$ws2.Range("A1").Value = "R"
$ws2.Range("B1").Value = "G"
$ws2.Range("C1").Value = "B"

# --- RGB sample rows (rows 2-33) ---
$data = @(
  @(97, 87, 75),
  @(115, 101, 98),
  @(119, 102, 95),
  @(126, 103, 97),
  @(196, 175, 158),
  @(159, 131, 127),
  @(190, 157, 150),
  @(191, 168, 154),
  @(223, 196, 179),
  @(184, 159, 139),
  @(133, 107, 90),
  @(111, 92, 75),
  @(127, 110, 100),
  @(112, 97, 90),
  @(100, 93, 85),
  @(89, 83, 83),
  @(93, 82, 96),
  @(99, 89, 88),
  @(98, 86, 86),
  @(114, 95, 88),
  @(126, 106, 97),
  @(91, 92, 97),
  @(121, 115, 115),
  @(103, 100, 93),
  @(133, 115, 113),
  @(122, 108, 97),
  @(123, 105, 105),
  @(209, 182, 175),
  @(177, 148, 130),
  @(114, 100, 99),
  @(124, 111, 103),
  @(110, 104, 108)
)

$r = 2
foreach ($row in $data) {
  $ws2.Cells.Item($r, 1).Value = $row[0]
  $ws2.Cells.Item($r, 2).Value = $row[1]
  $ws2.Cells.Item($r, 3).Value = $row[2]
  $r++
}

# --- Minimum-value summary row (row 34), shared formula across A:B, separate for C ---
$ws2.Range("A34:B34").Formula = "=MIN(A2:A33)"
$ws2.Range("C34").Formula = "=MIN(C2:C33)"
$ws2.Range("D34").Value = "minimum"

# --- Sheet2 view state: selection on A34:C34 ---
$ws2.Range("A34:C34").Select()

# --- Sheet1 ("frame3") selection moves from D43 to A42; re-activate it as the visible tab ---
$ws1.Activate()
$ws1.Range("A42").Select()
